$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fill in the two new test-case rows (34 = AuthID_0, 35 = AuthID_1) that were
# previously blank placeholder rows (already styled/merged/conditionally
# formatted, just waiting for data).
# ---------------------------------------------------------------------------

# -- Row 34 : AuthID_0 ------------------------------------------------------
$ws.Range("A34").Value = "Validate the `nfunctionality of `nauthenticating`nan ID"
$ws.Range("B34").Value = "AuthID_0"
$ws.Range("C34").Value = "Test if we can `nauthenticate a valid `nID"
$ws.Range("E34").Value = "1.STMCUBE_IDE`n2.Proteus Simulation`n3.GPIO MCAL driver`n4.LCD HAL driver`n5.Keypad HAL driver`n6.SYSTICK MCAL driver"
$ws.Range("F34").Value = "ID =`"000`"`nIDS_LIST = {`"000`"}"
$ws.Range("G34").Value = "Send a valid ID signal"
$ws.Range("H34").Value = "As expected"
$ws.Range("I34").Value = "Pass"
$ws.Range("J34").Value = "Hossam Eid"
$ws.Range("K34").Value = "Hossam Eid"
$ws.Range("L34").Value = "Functional Test"

# -- Row 35 : AuthID_1 (A is merged with A34, left blank) -------------------
$ws.Range("B35").Value = "AuthID_1"
$ws.Range("C35").Value = "Test if we can `nreject an invalid `nID"
$ws.Range("E35").Value = "1.STMCUBE_IDE`n2.Proteus Simulation`n3.GPIO MCAL driver`n4.LCD HAL driver`n5.Keypad HAL driver`n6.SYSTICK MCAL driver"
$ws.Range("F35").Value = "ID =`"111`"`nIDS_LIST = {`"000`"}"
$ws.Range("G35").Value = "Send an invalid ID`n signal"
$ws.Range("H35").Value = "As expected"
$ws.Range("I35").Value = "Pass"
$ws.Range("J35").Value = "Hossam Eid"
$ws.Range("K35").Value = "Hossam Eid"
$ws.Range("L35").Value = "Functional Test"

# ---------------------------------------------------------------------------
# Re-apply the formatting that the (already-templated) rows 34/35 need by
# copying it across from other rows in the sheet that already carry the
# exact same look.
# ---------------------------------------------------------------------------

# B/F/G/H/I/J/K/L columns match row 32 (the other single-row test group).
$ws.Range("B32").Copy()
$ws.Range("B34,B35").PasteSpecial(-4122)

$ws.Range("F32").Copy()
$ws.Range("F34,F35").PasteSpecial(-4122)

$ws.Range("G32").Copy()
$ws.Range("G34").PasteSpecial(-4122)

$ws.Range("H32").Copy()
$ws.Range("H34,H35").PasteSpecial(-4122)

$ws.Range("I32").Copy()
$ws.Range("I34,I35").PasteSpecial(-4122)

$ws.Range("J32").Copy()
$ws.Range("J34,J35").PasteSpecial(-4122)

$ws.Range("K32").Copy()
$ws.Range("K34,K35").PasteSpecial(-4122)

$ws.Range("L32").Copy()
$ws.Range("L34,L35").PasteSpecial(-4122)

$ws.Range("E32").Copy()
$ws.Range("E35").PasteSpecial(-4122)

# C column (center/vcenter/wrap, no border) matches e.g. C7; also used for G35.
$ws.Range("C7").Copy()
$ws.Range("C34,C35,G35").PasteSpecial(-4122)

# D column is an always-empty styled spacer column (center/vcenter, no wrap),
# matching B7/B32's horizontal+vertical centering without wrap.
$ws.Range("B7").Copy()
$ws.Range("D34,D35").PasteSpecial(-4122)

# M35 should get the plain white-fill style used throughout columns D/M in
# the rows-25-30 block (e.g. D25).
$ws.Range("D25").Copy()
$ws.Range("M35").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row heights grew to fit the new wrapped text.
# ---------------------------------------------------------------------------
$ws.Rows.Item(34).RowHeight = 90.75
$ws.Rows.Item(35).RowHeight = 90

# ---------------------------------------------------------------------------
# Restore the view state (scroll position / zoom / selection) that was saved
# with the workbook.
# ---------------------------------------------------------------------------
$ws.Range("G39").Select()
$excel.ActiveWindow.Zoom = 80
